# Updates cryptos list: Price (column D) and Volume(1h) (column E) cells
# for rows 2-51 on the active worksheet, per upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '45.800.57'
$ws.Cells.Item(2, 5).Value = '  +2.80%  '
$ws.Cells.Item(3, 4).Value = '2.439.38'
$ws.Cells.Item(3, 5).Value = '  +0.33%  '
$ws.Cells.Item(4, 5).Value = '  -0.01%  '
$ws.Cells.Item(5, 4).Value = '321.49'
$ws.Cells.Item(5, 5).Value = '  +2.98%  '
$ws.Cells.Item(6, 4).Value = '104.39'
$ws.Cells.Item(6, 5).Value = '  +2.34%  '
$ws.Cells.Item(7, 5).Value = '  +0.81%  '
$ws.Cells.Item(8, 4).Value = '0.999'
$ws.Cells.Item(8, 5).Value = '  -0.05%  '
$ws.Cells.Item(9, 4).Value = '0.538'
$ws.Cells.Item(9, 5).Value = '  +5.81%  '
$ws.Cells.Item(10, 4).Value = '35.96'
$ws.Cells.Item(11, 4).Value = '0.0805'
$ws.Cells.Item(11, 5).Value = '  +0.56%  '
$ws.Cells.Item(12, 5).Value = '  -1.12%  '
$ws.Cells.Item(13, 4).Value = '18.32'
$ws.Cells.Item(13, 5).Value = '  -2.50%  '
$ws.Cells.Item(14, 4).Value = '7.05'
$ws.Cells.Item(14, 5).Value = '  +1.00%  '
$ws.Cells.Item(15, 4).Value = '2.819.08'
$ws.Cells.Item(15, 5).Value = '  +0.25%  '
$ws.Cells.Item(16, 4).Value = '2.431.96'
$ws.Cells.Item(16, 5).Value = '  +0.26%  '
$ws.Cells.Item(17, 4).Value = '0.841'
$ws.Cells.Item(17, 5).Value = '  +0.14%  '
$ws.Cells.Item(18, 4).Value = '45.652.13'
$ws.Cells.Item(18, 5).Value = '  +2.65%  '
$ws.Cells.Item(19, 4).Value = '12.39'
$ws.Cells.Item(19, 5).Value = '  -0.71%  '
$ws.Cells.Item(20, 4).Value = '6.43'
$ws.Cells.Item(20, 5).Value = '  +0.14%  '
$ws.Cells.Item(21, 4).Value = '0.0₃0930'
$ws.Cells.Item(21, 5).Value = '  +2.22%  '
$ws.Cells.Item(22, 4).Value = '71.63'
$ws.Cells.Item(22, 5).Value = '  +3.93%  '
$ws.Cells.Item(23, 4).Value = '2.36'
$ws.Cells.Item(23, 5).Value = '  +1.48%  '
$ws.Cells.Item(24, 4).Value = '246.79'
$ws.Cells.Item(24, 5).Value = '  +2.49%  '
$ws.Cells.Item(25, 4).Value = '2.51'
$ws.Cells.Item(25, 5).Value = '  +0.18%  '
$ws.Cells.Item(26, 5).Value = '  +0.04%  '
$ws.Cells.Item(27, 4).Value = '25.79'
$ws.Cells.Item(27, 5).Value = '  +2.24%  '
$ws.Cells.Item(28, 4).Value = '2.19'
$ws.Cells.Item(28, 5).Value = '  -4.70%  '
$ws.Cells.Item(29, 4).Value = '9.67'
$ws.Cells.Item(29, 5).Value = '  -0.18%  '
$ws.Cells.Item(30, 4).Value = '33.47'
$ws.Cells.Item(30, 5).Value = '  +0.83%  '
$ws.Cells.Item(31, 4).Value = '49.35'
$ws.Cells.Item(31, 5).Value = '  +1.55%  '
$ws.Cells.Item(32, 5).Value = '  +5.02%  '
$ws.Cells.Item(33, 4).Value = '20.25'
$ws.Cells.Item(33, 5).Value = '  +3.45%  '
$ws.Cells.Item(34, 4).Value = '5.26'
$ws.Cells.Item(34, 5).Value = '  +1.06%  '
$ws.Cells.Item(35, 5).Value = '  -0.11%  '
$ws.Cells.Item(36, 5).Value = '  -0.89%  '
$ws.Cells.Item(37, 5).Value = '  +1.00%  '
$ws.Cells.Item(38, 5).Value = '  -0.75%  '
$ws.Cells.Item(39, 4).Value = '2.92'
$ws.Cells.Item(39, 5).Value = '  +0.54%  '
$ws.Cells.Item(40, 4).Value = '127.07'
$ws.Cells.Item(40, 5).Value = '  +0.40%  '
$ws.Cells.Item(41, 5).Value = '  -2.22%  '
$ws.Cells.Item(42, 5).Value = '  +1.34%  '
$ws.Cells.Item(43, 4).Value = '21.01'
$ws.Cells.Item(43, 5).Value = '  -4.36%  '
$ws.Cells.Item(44, 4).Value = '0.0292'
$ws.Cells.Item(44, 5).Value = '  +0.15%  '
$ws.Cells.Item(45, 4).Value = '1.960.58'
$ws.Cells.Item(45, 5).Value = '  +0.63%  '
$ws.Cells.Item(46, 5).Value = '  +0.79%  '
$ws.Cells.Item(47, 5).Value = '  -2.80%  '
$ws.Cells.Item(48, 4).Value = '1.82'
$ws.Cells.Item(48, 5).Value = '  +8.89%  '
$ws.Cells.Item(49, 4).Value = '9.14'
$ws.Cells.Item(49, 5).Value = '  -4.57%  '
$ws.Cells.Item(50, 4).Value = '77.52'
$ws.Cells.Item(50, 5).Value = '  +4.88%  '
$ws.Cells.Item(51, 4).Value = '4.86'
$ws.Cells.Item(51, 5).Value = '  +4.91%  '
